$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new date column U with the same header style/format as the other date columns.
# Use a leading apostrophe so Excel stores the date-like text as a literal string
# instead of auto-converting it to a date serial number, then copy the header
# formatting (bold/border/alignment) from the neighboring header cell.
$ws.Range("U1").Formula = "'2025-06-23"
$ws.Range("T1").Copy()
$ws.Range("U1").PasteSpecial(-4122)

# Update totals for the new column and mark attendance as absent ("❌")
$ws.Range("S2").Value = 16
$ws.Range("U2").Value = "❌"

$ws.Range("S3").Value = 16
$ws.Range("U3").Value = "❌"
